# Add "Comments" field for reagent / reagent_batch (resolves #469)
# Adds a new column J with a "Comments" header (styled like the other
# bold headers) and one comment value per data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, bold like the rest of the header row
$ws.Range("J2").Value = "Comments"
$ws.Range("J2").Font.Bold = $true
$ws.Range("J2").WrapText = $false

# New data values, one per reagent / reagent_batch row
$ws.Range("J3").Value = "ORB test comment 1"
$ws.Range("J4").Value = "ORB test comment 2"
$ws.Range("J5").Value = "ORB test comment 3"
$ws.Range("J6").Value = "ORB test comment 4"

# Move the active selection to the new last cell, like Calc would after
# typing into J6
$ws.Range("J6").Select() | Out-Null
